$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6-9 (the MuSCs/Resolving-Mac target-cluster rows that are no longer present)
$ws.Rows("6:9").Delete()

# Row 2: ECs -> Ccl2/Cxcr3 -> Resolving-Mac
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 4.232924
$ws.Range("H2").Value = 12.698772
$ws.Range("I2").Value = 0.05792409824508498
$ws.Range("J2").Value = 0.05792409824508497
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.888791333333333
$ws.Range("N2").Value = 5.666374
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 7.995110165858667
$ws.Range("R2").Value = 71.955991492728
$ws.Range("S2").Value = 0.05792409824508498
$ws.Range("T2").Value = 0.05792409824508497

# Row 3: FAPs -> Ccl2/Cxcr3 -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 23.77965533333333
$ws.Range("H3").Value = 71.338966
$ws.Range("I3").Value = 0.3254051080913003
$ws.Range("J3").Value = 0.3254051080913002
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.888791333333333
$ws.Range("N3").Value = 5.666374
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 44.91480690325378
$ws.Range("R3").Value = 404.233262129284
$ws.Range("S3").Value = 0.3254051080913003
$ws.Range("T3").Value = 0.3254051080913002

# Row 4: MuSCs -> Ccl2/Cxcr3 -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 7.006365333333332
$ws.Range("H4").Value = 21.019096
$ws.Range("I4").Value = 0.09587637148905993
$ws.Range("J4").Value = 0.09587637148905992
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.888791333333333
$ws.Range("N4").Value = 5.666374
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 13.23356211976711
$ws.Range("R4").Value = 119.102059077904
$ws.Range("S4").Value = 0.09587637148905993
$ws.Range("T4").Value = 0.09587637148905992

# Row 5: Resolving-Mac -> Ccl2/Cxcr3 -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 38.058136
$ws.Range("H5").Value = 114.174408
$ws.Range("I5").Value = 0.5207944221745548
$ws.Range("J5").Value = 0.5207944221745548
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.888791333333333
$ws.Range("N5").Value = 5.666374
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 71.88387743962133
$ws.Range("R5").Value = 646.9548969565921
$ws.Range("S5").Value = 0.5207944221745548
$ws.Range("T5").Value = 0.5207944221745548

Write-Output "done"
